$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 192, shifting existing rows 192:220 down to 193:221
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new data record
$ws.Range("A192").Value = 10
$ws.Range("B192").Value = "Vega Modelo de Temuco"
$ws.Range("C192").Value = "La Araucanía"
$ws.Range("D192").Value = 44505
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = 100112044
$ws.Range("G192").Value = "Perejil"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 45
$ws.Range("K192").Value = 5000
$ws.Range("L192").Value = 5000
$ws.Range("M192").Value = 5000
$ws.Range("N192").Value = "$/docena de atados (3 kilos)"
$ws.Range("O192").Value = "Provincia de Cautín"
$ws.Range("P192").Value = 1667
$ws.Range("Q192").Value = 3
$ws.Range("R192").Value = "Hortaliza"
